# Daily attendance processing - 2025-11-08 12:35:33
# Normalize the "Recorded By" (column G) values: the automated "System" /
# "admin@admin.com" marker that used to be listed first is moved so the
# human/service recorder's address is listed first instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Map of exact "before" strings to their corrected "after" strings, as used
# in column G ("Recorded By") throughout the sheet.
$replacements = @{
    "System, backup@backdoor.com, system" = "backup@backdoor.com, System, system"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $current = $cell.Value2

    if ($current -ne $null -and $replacements.ContainsKey($current)) {
        $cell.Value2 = $replacements[$current]
    }
}
